$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Add the new "Debts" worksheet at the end of the workbook
# ------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$debts = $wb.Worksheets.Add($null, $lastSheet)
$debts.Name = "Debts"

# Header row
$debts.Range("A1").Value = "name"
$debts.Range("B1").Value = "type"
$debts.Range("C1").Value = "principal"
$debts.Range("D1").Value = "paid"
$debts.Range("E1").Value = "creditor"
$debts.Range("F1").Value = "dueDate"
$debts.Range("G1").Value = "description"

# Data row
$debts.Range("A2").Value = "HOme"
$debts.Range("B2").Value = "Loan"
$debts.Range("C2").Value = 50000
$debts.Range("D2").Value = 0

# E2 / G2 are stored as empty strings (present but blank) rather than a
# truly empty/missing cell, so use an existing blank-string cell as the
# source to stamp that "blank text" cell type onto them.
$debts.Range("A1").Copy($debts.Range("X1"))
$debts.Range("X1").Value = ""
$debts.Range("X1").Copy($debts.Range("E2"))
$debts.Range("X1").Copy($debts.Range("G2"))
$debts.Range("X1").ClearContents()

# dueDate is kept as literal text "2026-01-02", not converted to a date
# serial number, so force the cell to Text format before assigning it.
$debts.Range("F2").NumberFormat = "@"
$debts.Range("F2").Value = "2026-01-02"

# ------------------------------------------------------------------
# 2. Update the "Expenses" sheet (sheet2.xml): row 3 gains blank
#    placeholders in columns F (expenseFor) and I (contributions),
#    matching the other blank text cells already present in that row.
# ------------------------------------------------------------------
$expenses = $wb.Worksheets.Item("Expenses")
$expenses.Range("J3").Copy($expenses.Range("F3"))
$expenses.Range("I2").Copy($expenses.Range("I3"))
